$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 34.71251733333333
$ws.Range("H2").Value = 104.137552
$ws.Range("I2").Value = 0.111750244749681
$ws.Range("J2").Value = 0.111750244749681
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 66.61006935687644
$ws.Range("R2").Value = 599.490624211888
$ws.Range("S2").Value = 0.00072906073816197
$ws.Range("T2").Value = 0.00072906073816197
$ws.Range("G3").Value = 34.71251733333333
$ws.Range("H3").Value = 104.137552
$ws.Range("I3").Value = 0.111750244749681
$ws.Range("J3").Value = 0.111750244749681
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 6292.97642649045
$ws.Range("R3").Value = 56636.78783841405
$ws.Range("S3").Value = 0.06887790514302129
$ws.Range("T3").Value = 0.06887790514302129
$ws.Range("G4").Value = 34.71251733333333
$ws.Range("H4").Value = 104.137552
$ws.Range("I4").Value = 0.111750244749681
$ws.Range("J4").Value = 0.111750244749681
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 1008.31782936382
$ws.Range("R4").Value = 9074.860464274385
$ws.Range("S4").Value = 0.01103624344000135
$ws.Range("T4").Value = 0.01103624344000135
$ws.Range("G5").Value = 34.71251733333333
$ws.Range("H5").Value = 104.137552
$ws.Range("I5").Value = 0.111750244749681
$ws.Range("J5").Value = 0.111750244749681
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 2842.070185541425
$ws.Range("R5").Value = 25578.63166987283
$ws.Range("S5").Value = 0.03110703542849642
$ws.Range("T5").Value = 0.03110703542849642
$ws.Range("H6").Value = 578.4917909999999
$ws.Range("I6").Value = 0.620780861354714
$ws.Range("J6").Value = 0.6207808613547139
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 370.0238538437476
$ws.Range("R6").Value = 3330.214684593728
$ws.Range("S6").Value = 0.004049986235196887
$ws.Range("T6").Value = 0.004049986235196886
$ws.Range("H7").Value = 578.4917909999999
$ws.Range("I7").Value = 0.620780861354714
$ws.Range("J7").Value = 0.6207808613547139
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("Q7").Value = 34957.94872997629
$ws.Range("S7").Value = 0.3826218490954588
$ws.Range("T7").Value = 0.3826218490954587
$ws.Range("H8").Value = 578.4917909999999
$ws.Range("I8").Value = 0.620780861354714
$ws.Range("J8").Value = 0.6207808613547139
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 5601.279997497049
$ws.Range("R8").Value = 50411.51997747344
$ws.Range("S8").Value = 0.06130714723847532
$ws.Range("T8").Value = 0.06130714723847532
$ws.Range("H9").Value = 578.4917909999999
$ws.Range("I9").Value = 0.620780861354714
$ws.Range("J9").Value = 0.6207808613547139
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 15787.90974250635
$ws.Range("R9").Value = 142091.1876825571
$ws.Range("S9").Value = 0.172801878785583
$ws.Range("T9").Value = 0.172801878785583
$ws.Range("G10").Value = 19.96051866666667
$ws.Range("H10").Value = 59.881556
$ws.Range("I10").Value = 0.06425903442584988
$ws.Range("J10").Value = 0.06425903442584986
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 38.30236568608489
$ws.Range("R10").Value = 344.721291174764
$ws.Range("S10").Value = 0.0004192271719585588
$ws.Range("T10").Value = 0.0004192271719585587
$ws.Range("G11").Value = 19.96051866666667
$ws.Range("H11").Value = 59.881556
$ws.Range("I11").Value = 0.06425903442584988
$ws.Range("J11").Value = 0.06425903442584986
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 3618.610319258972
$ws.Range("R11").Value = 32567.49287333075
$ws.Range("S11").Value = 0.03960642491369989
$ws.Range("T11").Value = 0.03960642491369988
$ws.Range("G12").Value = 19.96051866666667
$ws.Range("H12").Value = 59.881556
$ws.Range("I12").Value = 0.06425903442584988
$ws.Range("J12").Value = 0.06425903442584986
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 579.806605832717
$ws.Range("R12").Value = 5218.259452494452
$ws.Range("S12").Value = 0.006346101064312261
$ws.Range("T12").Value = 0.00634610106431226
$ws.Range("G13").Value = 19.96051866666667
$ws.Range("H13").Value = 59.881556
$ws.Range("I13").Value = 0.06425903442584988
$ws.Range("J13").Value = 0.06425903442584986
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 1634.2575920301
$ws.Range("R13").Value = 14708.3183282709
$ws.Range("S13").Value = 0.01788728127587916
$ws.Range("T13").Value = 0.01788728127587916
$ws.Range("G14").Value = 63.12224
$ws.Range("H14").Value = 189.36672
$ws.Range("I14").Value = 0.2032098594697551
$ws.Range("J14").Value = 0.2032098594697551
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 121.1256661101867
$ws.Range("R14").Value = 1090.13099499168
$ws.Range("S14").Value = 0.001325745017191408
$ws.Range("T14").Value = 0.001325745017191408
$ws.Range("G15").Value = 63.12224
$ws.Range("H15").Value = 189.36672
$ws.Range("I15").Value = 0.2032098594697551
$ws.Range("J15").Value = 0.2032098594697551
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 11443.32934695659
$ws.Range("R15").Value = 102989.9641226093
$ws.Range("S15").Value = 0.1252495639364086
$ws.Range("T15").Value = 0.1252495639364086
$ws.Range("G16").Value = 63.12224
$ws.Range("H16").Value = 189.36672
$ws.Range("I16").Value = 0.2032098594697551
$ws.Range("J16").Value = 0.2032098594697551
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 1833.554144466027
$ws.Range("R16").Value = 16501.98730019424
$ws.Range("S16").Value = 0.02006862252105342
$ws.Range("T16").Value = 0.02006862252105342
$ws.Range("G17").Value = 63.12224
$ws.Range("H17").Value = 189.36672
$ws.Range("I17").Value = 0.2032098594697551
$ws.Range("J17").Value = 0.2032098594697551
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 5168.102175531946
$ws.Range("R17").Value = 46512.91957978751
$ws.Range("S17").Value = 0.05656592799510173
$ws.Range("T17").Value = 0.05656592799510172
